$d = $word.ActiveDocument

# 1) Update the date in paragraph 1 (in-place substring replace)
$d.Content.Find.Execute("17.07.24", $false, $false, $false, $false, $false, $true, 1, $false, "16.07.24", 2) | Out-Null

# 2-6) Replace whole paragraph texts (exclude trailing paragraph mark)
$p = $d.Paragraphs.Item(2)
$r = $p.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = 'How Does Quantization Affect Multilingual LLMs?'

$p = $d.Paragraphs.Item(3)
$r = $p.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = 'היום נסקור קצרות מאמר שחוקר נושא חשוב לכל מי שעוסק במודלי שפה. הנושא הזה הוא קוונטיזציה או קווינטוט של מודלי שפה שמאפשר לנו גם להקטין את כמות הזכרון הנדרש לאחסון של המודל וגם מזרז את האינפרנס של המודל. אבל כמובן שזה לא בא בלי המחיר והמחיר הוא ביצועיי המודל. המאמר חוקר עד כמה חמורה פגיעה בביצועי המודלי לכמה רמות ושיטות קווינטוט(ניתן לקוונטט שכבות שונות ברמות שונות וגם ניתן לקוונטט משקלי המודל והאקטיבציות ברמות שונות של קווינטוט).'

$p = $d.Paragraphs.Item(4)
$r = $p.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = 'המאמר נכתב על ידי מדעני חברת cohere ובאופן טבעי מתמקד במודלי שלהם. המחברים לקחו מודלים בגדלים שונים ובדקו אותם במספר בנצ''מארקים שונים וגם ביצעו אבלואציה של ביצועי המודלים על ידי בודקים אנושיים. המחברים הגיעו למספר מסקנות מעניינות:'

$p = $d.Paragraphs.Item(5)
$r = $p.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = 'הפגיעה מהקווינטוט הנמדדת על הבנצ''מארקים משמעותית קטנה יותר מזו הנעשית על ידי בודקים אנושיים.'

$p = $d.Paragraphs.Item(6)
$r = $p.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = 'הפגיעה לרוב מחמירה ככל שקווינטוט נהיה יותר קשוח כלומר לפחות ביטים'

# 7-11) Append new paragraphs at the end, matching paragraph 6 style (Normal)
$n = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($n)
$r = $last.Range
$r.Collapse(0) | Out-Null
$r.InsertParagraphAfter() | Out-Null
$r.Collapse(0) | Out-Null
$r.InsertAfter('מודלים גדולים בד״כ עמידים יותר לקווינטוט מאשר מודלים קטנים יותר') | Out-Null

$n = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($n)
$r = $last.Range
$r.Collapse(0) | Out-Null
$r.InsertParagraphAfter() | Out-Null
$r.Collapse(0) | Out-Null
$r.InsertAfter('מודלים מולטי-שפתיים (multilingual) סובלים יותר מקווינטוט מאשר מודלים חד שפתיים והביצועים על השפות הפחות נפוצות נפגעות יותר מאשר על שפות נפוצות יותר') | Out-Null

$n = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($n)
$r = $last.Range
$r.Collapse(0) | Out-Null
$r.InsertParagraphAfter() | Out-Null
$r.Collapse(0) | Out-Null
$r.InsertAfter('היכולת של המודלי ל-reasoning (למשל יכולת לפתור שאלות מתמטיות) נפגעת מאוד מהקוויטוט.') | Out-Null

$n = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($n)
$r = $last.Range
$r.Collapse(0) | Out-Null
$r.InsertParagraphAfter() | Out-Null
$r.Collapse(0) | Out-Null
$r.InsertAfter('יש עוד כמה מציאות מעניינות…') | Out-Null

$n = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($n)
$r = $last.Range
$r.Collapse(0) | Out-Null
$r.InsertParagraphAfter() | Out-Null
$r.Collapse(0) | Out-Null
$r.InsertAfter('https://arxiv.org/pdf/2407.03211') | Out-Null

"OK paragraphs=" + $d.Paragraphs.Count